$wb = $excel.ActiveWorkbook

$wsFolha = $wb.Worksheets.Item("Folha Ponto")
$wsTotais = $wb.Worksheets.Item("Totais")

# --- Folha Ponto: fill in row 27 with the new time entry ---
$wsFolha.Range("B27").Value = 44906
$wsFolha.Range("C27").Value = 0.79166666666666663
$wsFolha.Range("D27").Value = 0.84375
$wsFolha.Range("F27").Value = "SITS"
$wsFolha.Range("G27").Value = "Organizando fases e testando"

# --- View/selection changes ---
# Totais becomes the active/selected sheet, Folha Ponto loses tabSelected.
$wsFolha.Activate()
$wsFolha.Range("D28").Select()

$wsTotais.Activate()
$wsTotais.Range("H2").Select()
